# Add a new solicitação row (row 4) to the "Solicitações" sheet,
# mirroring the existing rows' layout and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 4

# Date value (serial 45439 == 2024-05-27), formatted like A2/A3 (yyyy-mm-dd)
$ws.Cells.Item($row, 1).Value = 45439
$ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"

$ws.Cells.Item($row, 2).Value  = "Thayna Silva Santana"
$ws.Cells.Item($row, 3).Value  = "Rua Castro Alves"
$ws.Cells.Item($row, 4).Value  = "thayna.silva"
$ws.Cells.Item($row, 5).Value  = "Francisco"
$ws.Cells.Item($row, 6).Value  = "Estagiario"
$ws.Cells.Item($row, 7).Value  = "Automação"
$ws.Cells.Item($row, 8).Value  = "Todos"
$ws.Cells.Item($row, 9).Value  = "Notebook + Carregador"
$ws.Cells.Item($row, 10).Value = "LUM-001-001-012"
$ws.Cells.Item($row, 11).Value = "Mayara Almeilda"
$ws.Cells.Item($row, 12).Value = "23/05/2024 16:20:15"
